$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1238682.5
$ws.Range("I86").Value = 1589734.6
$ws.Range("K86").Value = 1589734.6
$ws.Range("M86").Value = -1588611.6
$ws.Range("H89").Value = 1238682.5
$ws.Range("I89").Value = 1589734.6
$ws.Range("K89").Value = 7948673
$ws.Range("M89").Value = -7943057
$ws.Range("H92").Value = 848.625
$ws.Range("I92").Value = 827.0833
$ws.Range("J92").Value = 913.25
$ws.Range("K92").Value = 827.0833
$ws.Range("L92").Value = 913.25
$ws.Range("M92").Value = 420.9167
$ws.Range("N92").Value = -3409.25
$ws.Range("H113").Value = 4500.25
$ws.Range("I113").Value = 4673
$ws.Range("J113").Value = 4080.7144
$ws.Range("K113").Value = 4673
$ws.Range("L113").Value = 4080.7144
$ws.Range("M113").Value = -1419
$ws.Range("N113").Value = -10588.7144
$ws.Range("H137").Value = 9547
$ws.Range("I137").Value = 5735.327
$ws.Range("J137").Value = 15553.272
$ws.Range("K137").Value = 17205.981
$ws.Range("L137").Value = 46659.81600000001
$ws.Range("M137").Value = -14655.981
$ws.Range("N137").Value = -51759.81600000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4159
$ws.Range("I2").Value = 3555.4167
$ws.Range("J2").Value = 5273.3076
$ws.Range("K2").Value = 3555.4167
$ws.Range("L2").Value = 5273.3076
$ws.Range("M2").Value = -3442.4167
$ws.Range("N2").Value = -5499.3076
$ws.Range("H45").Value = 8760.200000000001
$ws.Range("I45").Value = 10380.533
$ws.Range("J45").Value = 3899.2
$ws.Range("K45").Value = 10380.533
$ws.Range("L45").Value = 3899.2
$ws.Range("M45").Value = -10003.533
$ws.Range("N45").Value = -4653.2
$ws.Range("H61").Value = 3519.1853
$ws.Range("I61").Value = 1927.2273
$ws.Range("K61").Value = 1927.2273
$ws.Range("M61").Value = -1715.2273
$ws.Range("H97").Value = 1216.7273
$ws.Range("I97").Value = 1153.8889
$ws.Range("J97").Value = 1499.5
$ws.Range("K97").Value = 1153.8889
$ws.Range("L97").Value = 1499.5
$ws.Range("M97").Value = -657.8888999999999
$ws.Range("N97").Value = -2491.5
$ws.Range("H116").Value = 4159
$ws.Range("I116").Value = 3555.4167
$ws.Range("J116").Value = 5273.3076
$ws.Range("K116").Value = 3555.4167
$ws.Range("L116").Value = 5273.3076
$ws.Range("M116").Value = -1261.4167
$ws.Range("N116").Value = -9861.3076
$ws.Range("H132").Value = 10321.065
$ws.Range("J132").Value = 19978.143
$ws.Range("L132").Value = 59934.429
$ws.Range("N132").Value = -64994.429
$ws.Range("H136").Value = 3519.1853
$ws.Range("I136").Value = 1927.2273
$ws.Range("K136").Value = 5781.6819
$ws.Range("M136").Value = -3231.6819
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4159
$ws.Range("I3").Value = 3555.4167
$ws.Range("J3").Value = 5273.3076
$ws.Range("K3").Value = 3555.4167
$ws.Range("L3").Value = 5273.3076
$ws.Range("M3").Value = -3441.4167
$ws.Range("N3").Value = -5501.3076
$ws.Range("H22").Value = 629.8
$ws.Range("I22").Value = 587.5
$ws.Range("J22").Value = 799
$ws.Range("K22").Value = 587.5
$ws.Range("L22").Value = 799
$ws.Range("M22").Value = -414.5
$ws.Range("N22").Value = -1145
$ws.Range("H94").Value = 1058.7693
$ws.Range("I94").Value = 1086
$ws.Range("J94").Value = 1015.2
$ws.Range("K94").Value = 1086
$ws.Range("L94").Value = 1015.2
$ws.Range("M94").Value = -635
$ws.Range("N94").Value = -1917.2
$ws.Range("H134").Value = 7450.6055
$ws.Range("I134").Value = 5751.76
$ws.Range("J134").Value = 10717.615
$ws.Range("K134").Value = 17255.28
$ws.Range("L134").Value = 32152.845
$ws.Range("M134").Value = -14720.28
$ws.Range("N134").Value = -37222.845
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 273.59375
$ws.Range("I22").Value = 241.5
$ws.Range("J22").Value = 412.66666
$ws.Range("K22").Value = 241.5
$ws.Range("L22").Value = 412.66666
$ws.Range("M22").Value = 108.5
$ws.Range("N22").Value = -1112.66666
$ws.Range("H31").Value = 1659.6129
$ws.Range("I31").Value = 1427.1
$ws.Range("J31").Value = 1770.3334
$ws.Range("K31").Value = 1427.1
$ws.Range("L31").Value = 1770.3334
$ws.Range("M31").Value = -1132.1
$ws.Range("N31").Value = -2360.3334
$ws.Range("H34").Value = 1659.6129
$ws.Range("I34").Value = 1427.1
$ws.Range("J34").Value = 1770.3334
$ws.Range("K34").Value = 1427.1
$ws.Range("L34").Value = 1770.3334
$ws.Range("M34").Value = -1225.1
$ws.Range("N34").Value = -2174.3334
$ws.Range("H62").Value = 5449.4
$ws.Range("I62").Value = 6174
$ws.Range("J62").Value = 4966.3335
$ws.Range("K62").Value = 6174
$ws.Range("L62").Value = 4966.3335
$ws.Range("M62").Value = -5550
$ws.Range("N62").Value = -6214.3335
$ws.Range("H65").Value = 5449.4
$ws.Range("I65").Value = 6174
$ws.Range("J65").Value = 4966.3335
$ws.Range("K65").Value = 30870
$ws.Range("L65").Value = 24831.6675
$ws.Range("M65").Value = -27750
$ws.Range("N65").Value = -31071.6675
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37848508
$ws.Range("J4").Value = 116697.78
$ws.Range("L4").Value = 350093.34
$ws.Range("N4").Value = -350317.34
$ws.Range("H14").Value = 10656.571
$ws.Range("I14").Value = 10656.571
$ws.Range("K14").Value = 31969.713
$ws.Range("M14").Value = -31796.713
$ws.Range("H57").Value = 2999
$ws.Range("I57").Value = 2999
$ws.Range("K57").Value = 8997
$ws.Range("M57").Value = -8438
$ws.Range("H68").Value = 1836.9412
$ws.Range("J68").Value = 1758.2142
$ws.Range("L68").Value = 5274.642599999999
$ws.Range("N68").Value = -6896.642599999999
$ws.Range("H71").Value = 1836.9412
$ws.Range("J71").Value = 1758.2142
$ws.Range("L71").Value = 15823.9278
$ws.Range("N71").Value = -23935.9278
$ws.Range("H132").Value = 5002278
$ws.Range("I132").Value = 2365.8333
$ws.Range("J132").Value = 7145097.5
$ws.Range("K132").Value = 21292.4997
$ws.Range("L132").Value = 64305877.5
$ws.Range("M132").Value = -18762.4997
$ws.Range("N132").Value = -64310937.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2146.889
$ws.Range("I102").Value = 2040.375
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 2040.375
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = -418.375
$ws.Range("N102").Value = -6243
$ws.Range("H107").Value = 901.41174
$ws.Range("J107").Value = 1111.8572
$ws.Range("L107").Value = 1111.8572
$ws.Range("N107").Value = -4951.8572
$ws.Range("H122").Value = 1470.15
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2332.3333
$ws.Range("I126").Value = 2332.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6996.999899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4526.999899999999
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 3085
$ws.Range("I132").Value = 2948.3215
$ws.Range("K132").Value = 8844.9645
$ws.Range("M132").Value = -6314.9645
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2636.182
$ws.Range("I68").Value = 2424.8
$ws.Range("K68").Value = 2424.8
$ws.Range("M68").Value = -1675.8
$ws.Range("H71").Value = 2636.182
$ws.Range("I71").Value = 2424.8
$ws.Range("K71").Value = 12124
$ws.Range("M71").Value = -8380
$ws.Range("H82").Value = 1569.375
$ws.Range("I82").Value = 1061.6316
$ws.Range("K82").Value = 1061.6316
$ws.Range("M82").Value = -700.6315999999999
$ws.Range("H85").Value = 1569.375
$ws.Range("I85").Value = 1061.6316
$ws.Range("K85").Value = 1061.6316
$ws.Range("M85").Value = 186.3684000000001
$ws.Range("H136").Value = 27781400
$ws.Range("I136").Value = 55558052
$ws.Range("J136").Value = 4748.5
$ws.Range("K136").Value = 166674156
$ws.Range("L136").Value = 14245.5
$ws.Range("M136").Value = -166671606
$ws.Range("N136").Value = -19345.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 1659.5555
$ws.Range("I113").Value = 1659.5555
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4978.666499999999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2808.666499999999
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 4636.857
$ws.Range("I122").Value = 4454.154
$ws.Range("J122").Value = 5164.6665
$ws.Range("K122").Value = 13362.462
$ws.Range("L122").Value = 15493.9995
$ws.Range("M122").Value = -10912.462
$ws.Range("N122").Value = -20393.9995
$ws.Range("H126").Value = 3331.6743
$ws.Range("I126").Value = 3331.6743
$ws.Range("K126").Value = 9995.0229
$ws.Range("M126").Value = -7525.0229
$ws.Range("H132").Value = 13795.609
$ws.Range("I132").Value = 9775.692999999999
$ws.Range("J132").Value = 26257.35
$ws.Range("K132").Value = 29327.079
$ws.Range("L132").Value = 78772.04999999999
$ws.Range("M132").Value = -26797.079
$ws.Range("N132").Value = -83832.04999999999